$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New headers in row 5 (merged) ---
$ws.Range("H5").Value = "Battery = 8"
$ws.Range("H5:I5").Merge()
$ws.Range("J5").Value = "Battery = 10"
$ws.Range("J5:K5").Merge()
$ws.Range("H5:K5").HorizontalAlignment = -4108   # xlCenter
$ws.Range("H5:K5").Font.Bold = $true

# --- Duplicate J6 header (Battery) into J6/K6 as Average order / Number of vehicles ---
$ws.Range("J6").Value = "Average order"
$ws.Range("K6").Value = "Number of vehicles"
$ws.Range("J6:K6").HorizontalAlignment = -4108
$ws.Range("J6:K6").Font.Bold = $true

# --- Duplicate H:I values into J:K for rows 7-16 ---
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 4
$ws.Range("K8").Value = 5
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 4
$ws.Range("K10").Value = 5
$ws.Range("J11").Value = 60
$ws.Range("K11").Value = 5
$ws.Range("K12").Value = 6
$ws.Range("J13").Value = 70
$ws.Range("K13").Value = 6
$ws.Range("K14").Value = 7
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 6
$ws.Range("K16").Value = 7

$ws.Range("J7:K16").HorizontalAlignment = -4108

# --- New column G labels (merged) ---
$ws.Range("G7").Value = "Yuchen"
$ws.Range("G7:G10").Merge()
$ws.Range("G11").Value = "Hai"
$ws.Range("G11:G16").Merge()
$ws.Range("G7:G16").HorizontalAlignment = -4108

# --- Fill colors ---
$ws.Range("H7:K10").Interior.Color = 49087      # orange FFC000 (BGR 0x00C0FF -> decimal 49087)
$ws.Range("H11:K16").Interior.ThemeColor = 8
$ws.Range("H11:K16").Interior.TintAndShade = 0.59999389629810485

Write-Host "done"
